# Update cryptocurrency list values (price & 1h volume change) to reflect
# the latest scrape from coinranking.com, as performed by the scheduled
# GitHub Actions workflow.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.815.27'
$ws.Range('E2').Value = '  -0.59%  '
$ws.Range('D3').Value = '1.626.54'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('E4').Value = '  -0.30%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.92'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.64%  '
$ws.Range('E6').Value = '  -1.02%  '
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.17'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.41%  '
$ws.Range('E9').Value = '  -0.85%  '
$ws.Range('E10').Value = '  -1.16%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0878'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.39%  '
$ws.Range('D12').Value = '1.857.17'
$ws.Range('E12').Value = '  -0.45%  '
$ws.Range('D13').Value = '1.635.78'
$ws.Range('E13').Value = '  +0.09%  '
$ws.Range('E14').Value = '  -1.27%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.556'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.25%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.84'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.19%  '
$ws.Range('D17').Value = '27.836.70'
$ws.Range('E17').Value = '  -0.49%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '228.18'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.69%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.62'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.67%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0717'
$ws.Range('E20').Value = '  -1.16%  '
$ws.Range('E21').Value = '  -0.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.34'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.50%  '
$ws.Range('E23').Value = '  -4.95%  '
$ws.Range('E24').Value = '  -0.49%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '155.15'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.28%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.93'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.23%  '
$ws.Range('E27').Value = '  -0.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.43'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.39%  '
$ws.Range('E29').Value = '  -0.31%  '
$ws.Range('E30').Value = '  -0.64%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0480'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.50%  '
$ws.Range('E32').Value = '  -0.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.10'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.08%  '
$ws.Range('D34').Value = '1.406.35'
$ws.Range('E34').Value = '  -0.17%  '
$ws.Range('E35').Value = '  +2.17%  '
$ws.Range('E36').Value = '  -0.76%  '
$ws.Range('E37').Value = '  -1.48%  '
$ws.Range('E38').Value = '  -1.16%  '
$ws.Range('E39').Value = '  -1.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.848'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.995'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.23%  '
$ws.Range('E42').Value = '  -2.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '65.75'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.18%  '
$ws.Range('E44').Value = '  -0.65%  '
$ws.Range('E45').Value = '  -1.06%  '
$ws.Range('D46').Value = '1.766.43'
$ws.Range('E46').Value = '  -0.52%  '
$ws.Range('E47').Value = '  -3.87%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '88.34'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.34%  '
$ws.Range('E50').Value = '  -0.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.59'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.65%  '
